$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.651.26'
$ws.Range('E2').Value = '  +3.76%  '
$ws.Range('D3').Value = '3.504.91'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.503.65'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.97%  '
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('E11').Value = '  +4.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.440'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('D13').Value = '4.111.30'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').Value = '66.704.78'
$ws.Range('E17').Value = '  +3.80%  '
$ws.Range('D18').Value = '3.507.94'
$ws.Range('E18').Value = '  +4.57%  '
$ws.Range('E19').Value = '  +3.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '390.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.537'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000122'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.99%  '
$ws.Range('E27').Value = '  +7.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.181'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E31').Value = '  +6.03%  '
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.47'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.08%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +6.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.72'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('E38').Value = '  +3.08%  '
$ws.Range('E39').Value = '  +4.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0750'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.01%  '
$ws.Range('D44').Value = '2.811.73'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0312'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '353.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.60%  '
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.70'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.94%  '
